# Add a new "as of" forecast column (AB) and a new observed-date row (40)
# to both the "cases" and "deaths" worksheets, matching the daily sequence
# that the table already follows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# New diagonal / "as of" forecast values for column AB, keyed by row number,
# and the new bottom row (40) value, per sheet.
$abValues = @{
    "cases"  = @{ 27 = 84920; 28 = 92411; 29 = 100756; 30 = 108573; 31 = 118002; 32 = 127991; 33 = 138144; 34 = 147911; 35 = 157099; 36 = 165135; 37 = 172733; 38 = 180048; 39 = 186892; 40 = 193040 }
    "deaths" = @{ 27 = 5986;  28 = 6559;  29 = 7215;   30 = 7829;   31 = 8592;   32 = 9393;   33 = 10187;  34 = 10955;  35 = 11675;  36 = 12299;  37 = 12903;  38 = 13483;  39 = 14011;  40 = 14503  }
}

$b26Values = @{
    "cases"  = 78162
    "deaths" = 5466
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- New column header AB1: continues the daily date sequence as text ---
    $ws.Range("AB1").NumberFormat = "@"
    $ws.Range("AB1").Value = "2020-04-29"

    # --- Newly observed value for the existing row 26 ---
    $ws.Range("B26").Value = $b26Values[$name]

    # --- New forecast ("as of") values in column AB for rows 27-39 ---
    $vals = $abValues[$name]
    foreach ($r in 27..39) {
        $ws.Range("AB$r").Value = $vals[$r]
    }

    # --- New row 40 (new observed date) ---
    $ws.Range("A40").NumberFormat = "@"
    $ws.Range("A40").Value = "2020-05-13"
    $ws.Range("AB40").Value = $vals[40]

    # Normalize formatting / materialize blank cells for the new column by
    # copying the format from the neighboring column AA (rows 1-39), then
    # do the same for the brand new row 40 by copying row 39's formats.
    $ws.Range("AA1:AA39").Copy()
    $ws.Range("AB1:AB39").PasteSpecial(-4122)

    $ws.Range("A39:AB39").Copy()
    $ws.Range("A40:AB40").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
